$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "bruce@gmail.com"
$ws.Range("B2").Value = "Bruce"
$ws.Range("C2").Value = "Wayne"
$ws.Range("D2").Value = 45131.5644113271
$ws.Range("D3").Value = 45131.60717907707
$ws.Range("D4").Value = 45131.60718010938
$ws.Range("D5").Value = 45131.60718113895
